$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the transfer-purpose description (D8) with the new, longer report description.
$ws.Range("D8").Value = "Elaboración de un informe técnico sobre las capturas del tráfico de red del SmartHome, la imagen del disco del Raspberry Pi, el informe de diagnóstico de Google OnHub, los datos de Amazon Echo Alexa y  las adquisiciones de los dispositivos móviles de la víctima y de su marido."

# Update the transfer date (D10) to the new date.
$ws.Range("D10").Value = "Lunes 6 de Mayo de 2024"

# The longer description text requires a taller row to display properly.
$ws.Rows.Item(8).RowHeight = 120

# Move the selection to match where the author last left the cursor.
[void]$ws.Range("D12").Select()
